$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Hello!"
$ws.Range("A3").Interior.Color = 49407
[void]$ws.Range("A3").Select()
